# Update the PC1/PC2 attenuation model values in Sheet1.
# The underlying data generation was re-run with a more generalized variable
# (material properties in addition to crack geometry), producing new
# PC1 (column B) and PC2 (column C) values for each frequency row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B2").Value = -0.08794879144565823
$ws.Range("C2").Value = -0.2281345837858737
$ws.Range("B3").Value = -0.2689146979341474
$ws.Range("C3").Value = 0.2904403394249926
$ws.Range("B4").Value = -0.7173650107339983
$ws.Range("C4").Value = 0.2040834492862479
$ws.Range("B5").Value = -0.0990421880138692
$ws.Range("C5").Value = -0.3640027109257797
$ws.Range("B6").Value = 0.1031631439203668
$ws.Range("C6").Value = 0.2375362733434154
$ws.Range("B7").Value = -0.1495791038774646
$ws.Range("C7").Value = 0.3653367041317052
$ws.Range("B8").Value = -0.2189872995657428
$ws.Range("C8").Value = -0.3676778016270782
$ws.Range("B9").Value = -0.1169763045555941
$ws.Range("C9").Value = -0.2053920682659564
$ws.Range("B10").Value = 0.1468951921376749
$ws.Range("C10").Value = 0.300179750955637
$ws.Range("B11").Value = -0.2606168096992059
$ws.Range("C11").Value = 0.1146119600183457
$ws.Range("B12").Value = -0.1413993731241623
$ws.Range("C12").Value = -0.4233032765825469
$ws.Range("B13").Value = -0.1904091055034978
$ws.Range("C13").Value = 0.07622373906957672
$ws.Range("B14").Value = 0.2513379474984697
$ws.Range("C14").Value = 0.06768159240201649
$ws.Range("B15").Value = -0.1483007481057371
$ws.Range("C15").Value = -0.09069671435285163
$ws.Range("B16").Value = 0.2612963478696874
$ws.Range("C16").Value = -0.007641430201288569
$ws.Range("B17").Value = -0.04103851369585835
$ws.Range("C17").Value = -0.1451865937329306
